$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.8222690533928816
$ws.Range("J2").Value = 0.8222690533928814
$ws.Range("M2").Value = 0.01376766666666667
$ws.Range("N2").Value = 0.041303
$ws.Range("O2").Value = 0.02028832877083762
$ws.Range("P2").Value = 0.02028832877083762
$ws.Range("Q2").Value = 0.03175858802944444
$ws.Range("R2").Value = 0.285827292265
$ws.Range("S2").Value = 0.01668246489332021
$ws.Range("T2").Value = 0.01668246489332021

# Row 3
$ws.Range("I3").Value = 0.8222690533928816
$ws.Range("J3").Value = 0.8222690533928814
$ws.Range("M3").Value = 0.6648326666666667
$ws.Range("N3").Value = 1.994498
$ws.Range("O3").Value = 0.9797116712291625
$ws.Range("P3").Value = 0.9797116712291624
$ws.Range("Q3").Value = 1.533603861887778
$ws.Range("R3").Value = 13.80243475699
$ws.Range("S3").Value = 0.8055865884995614
$ws.Range("T3").Value = 0.8055865884995611

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4985973333333333
$ws.Range("H4").Value = 1.495792
$ws.Range("I4").Value = 0.1777309466071186
$ws.Range("J4").Value = 0.1777309466071185
$ws.Range("M4").Value = 0.01376766666666667
$ws.Range("N4").Value = 0.041303
$ws.Range("O4").Value = 0.02028832877083762
$ws.Range("P4").Value = 0.02028832877083762
$ws.Range("Q4").Value = 0.006864521886222222
$ws.Range("R4").Value = 0.061780696976
$ws.Range("S4").Value = 0.003605863877517409
$ws.Range("T4").Value = 0.003605863877517406

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4985973333333333
$ws.Range("H5").Value = 1.495792
$ws.Range("I5").Value = 0.1777309466071186
$ws.Range("J5").Value = 0.1777309466071185
$ws.Range("M5").Value = 0.6648326666666667
$ws.Range("N5").Value = 1.994498
$ws.Range("O5").Value = 0.9797116712291625
$ws.Range("P5").Value = 0.9797116712291624
$ws.Range("Q5").Value = 0.3314837947128889
$ws.Range("R5").Value = 2.983354152416
$ws.Range("S5").Value = 0.1741250827296012
$ws.Range("T5").Value = 0.1741250827296011
